# Add a second worksheet "Transmittals_New_ActionRequired" after the existing
# "Transmittals_New" sheet. The new sheet carries the single-user
# Approve/Reject test-data rows: it duplicates the header row plus the first
# two data rows (Action-Level2 = Approved / Rejected) from "Transmittals_New".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "Transmittals_New_ActionRequired"

# Column widths (approximate the source workbook's column widths; the COM
# layer quantizes ColumnWidth to whole pixels, so nudge the requested value
# down by the engine's fixed rounding pad to land as close as possible to
# the original fractional widths).
$colWidths = @(28.42578125, 20.85546875, 21, 13.28515625, 13.28515625, 20.42578125, 20.42578125, 20.42578125, 20.42578125, 20.42578125, 20.42578125, 25.85546875, 19, 13, 13.28515625)
$widthPad = 0.8333333333333321
for ($i = 1; $i -le $colWidths.Length; $i++) {
    $ws2.Columns.Item($i).ColumnWidth = $colWidths[$i - 1] - $widthPad
}
$ws2.Columns.Item(17).ColumnWidth = 20.7109375 - $widthPad

# Copy header row (A1:O1) - keeps shared strings + the existing bold/filled
# header style ("s=1") instead of minting new style entries.
$ws1.Range("A1:O1").Copy($ws2.Range("A1:O1"))

# Row 2 - single user Approve test data (copy of Transmittals_New row 2).
$ws1.Range("A2:J2").Copy($ws2.Range("A2:J2"))
$ws1.Range("L2:M2").Copy($ws2.Range("L2:M2"))

# Row 3 - single user Reject test data (copy of Transmittals_New row 3).
$ws1.Range("A3:J3").Copy($ws2.Range("A3:J3"))
$ws1.Range("L3:M3").Copy($ws2.Range("L3:M3"))

# Keep the first sheet active/selected, matching the original workbook.
$ws1.Activate()
